$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to remain a TEXT cell (matches original inlineStr/"t=s" cells),
    # so numeric-looking strings like "1.100" or "0.000008670" are not coerced into
    # doubles (which would silently drop significant trailing zeros / re-notate).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2
Set-TextValue "D2" '26.118.36'
Set-TextValue "E2" '  -3.42%  '

# Row 3
Set-TextValue "D3" '1.768.08'
Set-TextValue "E3" '  -1.60%  '

# Row 4
Set-TextValue "D4" '1.008'
Set-TextValue "E4" '  +0.32%  '

# Row 5
Set-TextValue "D5" '1.005'
Set-TextValue "E5" '  +0.11%  '

# Row 6
Set-TextValue "D6" '307.54'
Set-TextValue "E6" '  -0.12%  '

# Row 7
Set-TextValue "D7" '0.4213'
Set-TextValue "E7" '  +0.90%  '

# Row 8
Set-TextValue "D8" '0.3584'
Set-TextValue "E8" '  +0.98%  '

# Row 9
Set-TextValue "D9" '0.07103'
Set-TextValue "E9" '  +0.48%  '

# Row 10
Set-TextValue "D10" '0.8332'
Set-TextValue "E10" '  -1.43%  '

# Row 11
Set-TextValue "E11" '  +0.77%  '

# Row 12
Set-TextValue "D12" '1.743.09'
Set-TextValue "E12" '  -2.08%  '

# Row 13
Set-TextValue "D13" '6.423'
Set-TextValue "E13" '  +1.33%  '

# Row 14
Set-TextValue "D14" '5.209'
Set-TextValue "E14" '  -1.21%  '

# Row 15
Set-TextValue "D15" '0.06917'
Set-TextValue "E15" '  +2.44%  '

# Row 16
Set-TextValue "D16" '1.011'
Set-TextValue "E16" '  +0.36%  '

# Row 17
Set-TextValue "D17" '78.58'
Set-TextValue "E17" '  -1.58%  '

# Row 18
Set-TextValue "D18" '0.000008670'
Set-TextValue "E18" '  +0.04%  '

# Row 19
Set-TextValue "E19" '  -0.13%  '

# Row 20
Set-TextValue "D20" '14.88'
Set-TextValue "E20" '  -0.97%  '

# Row 21
Set-TextValue "D21" '26.123.35'
Set-TextValue "E21" '  -4.24%  '

# Row 22
Set-TextValue "D22" '5.086'
Set-TextValue "E22" '  +0.85%  '

# Row 23
Set-TextValue "D23" '10.93'
Set-TextValue "E23" '  -0.23%  '

# Row 24
Set-TextValue "D24" '1.980.84'
Set-TextValue "E24" '  -2.77%  '

# Row 25
Set-TextValue "D25" '151.65'
Set-TextValue "E25" '  -0.95%  '

# Row 26
Set-TextValue "D26" '1.802'
Set-TextValue "E26" '  -7.31%  '

# Row 27
Set-TextValue "D27" '17.94'
Set-TextValue "E27" '  -0.91%  '

# Row 28
Set-TextValue "D28" '5.029'
Set-TextValue "E28" '  +1.10%  '

# Row 29
Set-TextValue "D29" '113.90'
Set-TextValue "E29" '  +0.61%  '

# Row 30
Set-TextValue "D30" '1.843'
Set-TextValue "E30" '  +12.36%  '

# Row 31
Set-TextValue "D31" '0.08843'
Set-TextValue "E31" '  -0.69%  '

# Row 32
Set-TextValue "D32" '0.7198'
Set-TextValue "E32" '  +0.51%  '

# Row 33
Set-TextValue "D33" '1.111'
Set-TextValue "E33" '  +3.55%  '

# Row 34
Set-TextValue "D34" '4.285'
Set-TextValue "E34" '  -0.67%  '

# Row 35
Set-TextValue "D35" '1.006'
Set-TextValue "E35" '  +0.14%  '

# Row 36
Set-TextValue "D36" '2.749'
Set-TextValue "E36" '  -3.77%  '

# Row 37
Set-TextValue "D37" '1.100'
Set-TextValue "E37" '  +2.33%  '

# Row 38
Set-TextValue "D38" '0.05077'
Set-TextValue "E38" '  -0.57%  '

# Row 39
Set-TextValue "D39" '0.01873'
Set-TextValue "E39" '  -1.07%  '

# Row 40
Set-TextValue "D40" '0.1602'
Set-TextValue "E40" '  -0.93%  '

# Row 41
Set-TextValue "D41" '0.4887'
Set-TextValue "E41" '  -0.79%  '

# Row 42
Set-TextValue "D42" '2.583'
Set-TextValue "E42" '  +0.05%  '

# Row 43
Set-TextValue "D43" '6.308'
Set-TextValue "E43" '  +3.54%  '

# Row 44
Set-TextValue "D44" '7.975'
Set-TextValue "E44" '  -0.35%  '

# Row 45
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D45" '1.006'
Set-TextValue "E45" '  +0.18%  '

# Row 46
Set-TextValue "D46" '104.30'
Set-TextValue "E46" '  -0.12%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D47" '10.13'
Set-TextValue "E47" '  -0.08%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D48" '0.06167'
Set-TextValue "E48" '  -2.17%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D49" '1.610'
Set-TextValue "E49" '  +1.30%  '

# Row 50
Set-TextValue "D50" '0.4439'
Set-TextValue "E50" '  -1.25%  '

# Row 51
Set-TextValue "D51" '1.687'
Set-TextValue "E51" '  +0.90%  '
